# "Informe previo - simulacion 5"
# Update the subamortiguado input parameters (Ka, Cero, R6) and refresh
# which resistor row is highlighted in bold.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ka
$ws.Range("B1").Value = 1.245
# Cero
$ws.Range("E1").Value = 0.0000012
# R6
$ws.Range("E2").Value = 160

# Move the bold "chosen value" emphasis from row 12 (100uF) to row 9 (10uF)
$ws.Range("A9:C9").Font.Bold = $true
$ws.Range("A12:C12").Font.Bold = $false

# Leave the active cell on C9, matching where the user ended up
$ws.Range("C9").Select() | Out-Null
